$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 361, pushing the old
# rows 361-368 down to 363-370.
$ws.Rows("361:362").Insert()

# New row 361 (Primera) - week of 2021-09-09 (serial 44448)
$ws.Range("A361").Value = 8
$ws.Range("B361").Value = "Terminal La Palmera de La Serena"
$ws.Range("C361").Value = "Coquimbo"
$ws.Range("D361").Value = 44448
$ws.Range("E361").Value = 4
$ws.Range("F361").Value = 100112008
$ws.Range("G361").Value = "Coliflor"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 2300
$ws.Range("K361").Value = 650
$ws.Range("L361").Value = 700
$ws.Range("M361").Value = 675
$ws.Range("N361").Value = "$/unidad"
$ws.Range("O361").Value = "Provincia del Elquí"
$ws.Range("P361").Value = 675
$ws.Range("Q361").Value = 1
$ws.Range("R361").Value = "Hortaliza"

# New row 362 (Segunda) - same week, serial 44448
$ws.Range("A362").Value = 8
$ws.Range("B362").Value = "Terminal La Palmera de La Serena"
$ws.Range("C362").Value = "Coquimbo"
$ws.Range("D362").Value = 44448
$ws.Range("E362").Value = 4
$ws.Range("F362").Value = 100112008
$ws.Range("G362").Value = "Coliflor"
$ws.Range("H362").Value = "Sin especificar"
$ws.Range("I362").Value = "Segunda"
$ws.Range("J362").Value = 1340
$ws.Range("K362").Value = 550
$ws.Range("L362").Value = 600
$ws.Range("M362").Value = 575
$ws.Range("N362").Value = "$/unidad"
$ws.Range("O362").Value = "Provincia del Elquí"
$ws.Range("P362").Value = 575
$ws.Range("Q362").Value = 1
$ws.Range("R362").Value = "Hortaliza"
